$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.452.44'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.332.78'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'" + '305.02'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').Value = "'" + '101.56'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('D7').Value = "'" + '0.510'
$ws.Range('E7').Value = '  -2.79%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'" + '0.508'
$ws.Range('E9').Value = '  -2.67%  '
$ws.Range('D10').Value = "'" + '35.33'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = "'" + '6.80'
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').Value = '2.692.86'
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').Value = "'" + '15.64'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '2.346.24'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = "'" + '0.805'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '43.383.03'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = "'" + '6.11'
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('D22').Value = "'" + '68.27'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = "'" + '237.81'
$ws.Range('E23').Value = '  -1.66%  '
$ws.Range('E24').Value = '  -3.36%  '
$ws.Range('D25').Value = "'" + '2.55'
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = "'" + '25.02'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').Value = "'" + '2.29'
$ws.Range('E28').Value = '  +3.26%  '
$ws.Range('D29').Value = "'" + '34.63'
$ws.Range('E29').Value = '  -5.11%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'" + '9.23'
$ws.Range('E30').Value = '  -3.79%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'" + '165.08'
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').Value = "'" + '5.06'
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('D34').Value = "'" + '4.56'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('D36').Value = "'" + '17.08'
$ws.Range('E36').Value = '  -6.41%  '
$ws.Range('D37').Value = "'" + '0.0708'
$ws.Range('E37').Value = '  -4.17%  '
$ws.Range('D38').Value = "'" + '2.92'
$ws.Range('E38').Value = '  -6.68%  '
$ws.Range('E39').Value = '  -5.16%  '
$ws.Range('E40').Value = '  -3.87%  '
$ws.Range('E41').Value = '  -3.09%  '
$ws.Range('D42').Value = "'" + '2.48'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = '1.978.39'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('D45').Value = "'" + '18.70'
$ws.Range('E45').Value = '  -6.47%  '
$ws.Range('D46').Value = "'" + '10.04'
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('E47').Value = '  -4.89%  '
$ws.Range('D48').Value = "'" + '55.88'
$ws.Range('E48').Value = '  -4.97%  '
$ws.Range('D49').Value = "'" + '4.82'
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').Value = '2.555.37'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  -1.86%  '
